# Auto-generated-by-design edit script: applies the "Atualizacao de bases das ligas" diff
# for Romania Liga I.xlsx (rows 219-226 updated, rows 227-230 appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Update existing rows 219-226: fill in match results (H/I/J) and
#    refresh odds/closing-line columns (K..AC) to the latest snapshot.
# ------------------------------------------------------------------
# Row 219
$ws.Range("H219").Value = 1
$ws.Range("I219").Value = 1
$ws.Range("J219").Value = 'D'
$ws.Range("N219").Value = 1.95
$ws.Range("P219").Value = 3.1
$ws.Range("R219").Value = 2.05
$ws.Range("S219").Value = 1.8
$ws.Range("U219").Value = 1.925
$ws.Range("V219").Value = 1.925
$ws.Range("W219").Value = -1
$ws.Range("X219").Value = 2.5
$ws.Range("Y219").Value = -1
$ws.Range("Z219").Value = -1
$ws.Range("AA219").Value = 0.8
$ws.Range("AB219").Value = -1
$ws.Range("AC219").Value = 0.925

# Row 220
$ws.Range("H220").Value = 1
$ws.Range("I220").Value = 0
$ws.Range("J220").Value = 'H'
$ws.Range("N220").Value = 3.1
$ws.Range("O220").Value = 2.875
$ws.Range("P220").Value = 2.45
$ws.Range("Q220").Value = 0.25
$ws.Range("R220").Value = 1.75
$ws.Range("S220").Value = 2.125
$ws.Range("U220").Value = 2.05
$ws.Range("V220").Value = 1.8
$ws.Range("W220").Value = 2.1
$ws.Range("X220").Value = -1
$ws.Range("Y220").Value = -1
$ws.Range("Z220").Value = 0.75
$ws.Range("AA220").Value = -1
$ws.Range("AB220").Value = -1
$ws.Range("AC220").Value = 0.8

# Row 221
$ws.Range("H221").Value = 2
$ws.Range("I221").Value = 2
$ws.Range("J221").Value = 'D'
$ws.Range("N221").Value = 4.75
$ws.Range("O221").Value = 3.5
$ws.Range("P221").Value = 1.666
$ws.Range("R221").Value = 1.975
$ws.Range("S221").Value = 1.875
$ws.Range("U221").Value = 1.825
$ws.Range("V221").Value = 2.025
$ws.Range("W221").Value = -1
$ws.Range("X221").Value = 2.5
$ws.Range("Y221").Value = -1
$ws.Range("Z221").Value = 0.9750000000000001
$ws.Range("AA221").Value = -1
$ws.Range("AB221").Value = 0.825
$ws.Range("AC221").Value = -1

# Row 222
$ws.Range("H222").Value = 1
$ws.Range("I222").Value = 0
$ws.Range("J222").Value = 'H'
$ws.Range("N222").Value = 2.2
$ws.Range("R222").Value = 1.95
$ws.Range("S222").Value = 1.9
$ws.Range("U222").Value = 1.95
$ws.Range("V222").Value = 1.9
$ws.Range("W222").Value = 1.2
$ws.Range("X222").Value = -1
$ws.Range("Y222").Value = -1
$ws.Range("Z222").Value = 0.95
$ws.Range("AA222").Value = -1
$ws.Range("AB222").Value = -1
$ws.Range("AC222").Value = 0.8999999999999999

# Row 223
$ws.Range("H223").Value = 4
$ws.Range("I223").Value = 1
$ws.Range("J223").Value = 'H'
$ws.Range("N223").Value = 1.666
$ws.Range("P223").Value = 4.5
$ws.Range("R223").Value = 1.95
$ws.Range("S223").Value = 1.9
$ws.Range("T223").Value = 2.75
$ws.Range("U223").Value = 1.85
$ws.Range("V223").Value = 2
$ws.Range("W223").Value = 0.6659999999999999
$ws.Range("X223").Value = -1
$ws.Range("Y223").Value = -1
$ws.Range("Z223").Value = 0.95
$ws.Range("AA223").Value = -1
$ws.Range("AB223").Value = 0.8500000000000001
$ws.Range("AC223").Value = -1

# Row 224
$ws.Range("N224").Value = 3.75
$ws.Range("O224").Value = 2.9
$ws.Range("P224").Value = 2.05
$ws.Range("Q224").Value = 0.25
$ws.Range("R224").Value = 2.05
$ws.Range("S224").Value = 1.8
$ws.Range("U224").Value = 2.1
$ws.Range("V224").Value = 1.775

# Row 225
$ws.Range("N225").Value = 3.1
$ws.Range("O225").Value = 3.3
$ws.Range("P225").Value = 2.2
$ws.Range("R225").Value = 1.875
$ws.Range("S225").Value = 1.975
$ws.Range("U225").Value = 1.825
$ws.Range("V225").Value = 2.025

# Row 226
$ws.Range("N226").Value = 4.333
$ws.Range("O226").Value = 3.5
$ws.Range("P226").Value = 1.727
$ws.Range("Q226").Value = 0.75
$ws.Range("R226").Value = 1.875
$ws.Range("S226").Value = 1.975
$ws.Range("U226").Value = 1.875
$ws.Range("V226").Value = 1.975

# ------------------------------------------------------------------
# 2) Append 4 new fixtures as rows 227-230. Copy formats (bold/border
#    style on id column, date-time format on Date column) from the
#    last existing row (226) so the new rows match the sheet style
#    without introducing any new style entries.
# ------------------------------------------------------------------
$ws.Range("A226").Copy() | Out-Null
$ws.Range("A227:A230").PasteSpecial(-4122) | Out-Null
$ws.Range("E226").Copy() | Out-Null
$ws.Range("E227:E230").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 227
$ws.Range("A227").Value = 225
$ws.Range("B227").Value = 6836274
$ws.Range("C227").Value = 'Romania Liga I'
$ws.Range("D227").Value = 'Romania Liga I'
$ws.Range("E227").Value = 45352.625
$ws.Range("F227").Value = 'CSM Politehnica Iasi'
$ws.Range("G227").Value = 'Dinamo Bucharest'
$ws.Range("K227").Value = 2.4
$ws.Range("L227").Value = 3.1
$ws.Range("M227").Value = 2.875
$ws.Range("N227").Value = 2.3
$ws.Range("O227").Value = 3.1
$ws.Range("P227").Value = 3
$ws.Range("Q227").Value = -0.25
$ws.Range("R227").Value = 2.025
$ws.Range("S227").Value = 1.825
$ws.Range("T227").Value = 2.25
$ws.Range("U227").Value = 1.925
$ws.Range("V227").Value = 1.925
$ws.Range("W227").Value = 0
$ws.Range("X227").Value = 0
$ws.Range("Y227").Value = 0
$ws.Range("Z227").Value = 0
$ws.Range("AA227").Value = 0

# Row 228
$ws.Range("A228").Value = 226
$ws.Range("B228").Value = 6836275
$ws.Range("C228").Value = 'Romania Liga I'
$ws.Range("D228").Value = 'Romania Liga I'
$ws.Range("E228").Value = 45353.375
$ws.Range("F228").Value = 'ACS UTA Batrana Doamna'
$ws.Range("G228").Value = 'FC U Craiova 1948'
$ws.Range("K228").Value = 2.25
$ws.Range("L228").Value = 3.2
$ws.Range("M228").Value = 3
$ws.Range("N228").Value = 2.45
$ws.Range("O228").Value = 3.2
$ws.Range("P228").Value = 2.75
$ws.Range("Q228").Value = 0
$ws.Range("R228").Value = 1.8
$ws.Range("S228").Value = 2.05
$ws.Range("T228").Value = 2.25
$ws.Range("U228").Value = 1.925
$ws.Range("V228").Value = 1.925
$ws.Range("W228").Value = 0
$ws.Range("X228").Value = 0
$ws.Range("Y228").Value = 0
$ws.Range("Z228").Value = 0
$ws.Range("AA228").Value = 0

# Row 229
$ws.Range("A229").Value = 227
$ws.Range("B229").Value = 6836273
$ws.Range("C229").Value = 'Romania Liga I'
$ws.Range("D229").Value = 'Romania Liga I'
$ws.Range("E229").Value = 45353.5
$ws.Range("F229").Value = 'AFC Hermannstadt'
$ws.Range("G229").Value = 'Otelul Galati'
$ws.Range("K229").Value = 2.1
$ws.Range("L229").Value = 3
$ws.Range("M229").Value = 3.6
$ws.Range("N229").Value = 2.1
$ws.Range("O229").Value = 3
$ws.Range("P229").Value = 3.4
$ws.Range("Q229").Value = -0.25
$ws.Range("R229").Value = 1.85
$ws.Range("S229").Value = 2
$ws.Range("T229").Value = 2
$ws.Range("U229").Value = 1.925
$ws.Range("V229").Value = 1.925
$ws.Range("W229").Value = 0
$ws.Range("X229").Value = 0
$ws.Range("Y229").Value = 0
$ws.Range("Z229").Value = 0
$ws.Range("AA229").Value = 0

# Row 230
$ws.Range("A230").Value = 228
$ws.Range("B230").Value = 6907421
$ws.Range("C230").Value = 'Romania Liga I'
$ws.Range("D230").Value = 'Romania Liga I'
$ws.Range("E230").Value = 45353.625
$ws.Range("F230").Value = 'Farul Constanta'
$ws.Range("G230").Value = 'CFR Cluj'
$ws.Range("K230").Value = 2.7
$ws.Range("L230").Value = 3.2
$ws.Range("M230").Value = 2.5
$ws.Range("N230").Value = 2.75
$ws.Range("O230").Value = 3.2
$ws.Range("P230").Value = 2.45
$ws.Range("Q230").Value = 0
$ws.Range("R230").Value = 2.05
$ws.Range("S230").Value = 1.8
$ws.Range("T230").Value = 2.25
$ws.Range("U230").Value = 1.925
$ws.Range("V230").Value = 1.925
$ws.Range("W230").Value = 0
$ws.Range("X230").Value = 0
$ws.Range("Y230").Value = 0
$ws.Range("Z230").Value = 0
$ws.Range("AA230").Value = 0

Write-Output "edit applied"
